# Auto-generated edit script applying the cryptos.xlsx diff
# (cryptocurrency price/volume refresh + 3 row reorders, commit:
#  "Updated cryptos list on Tue Dec 19 01:14:06 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '42.691.72'
$ws.Cells.Item(2, 5).Value = '  +2.98%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.219.19'
$ws.Cells.Item(3, 5).Value = '  +1.00%  '
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '244.09'
$ws.Cells.Item(5, 5).Value = '  +1.88%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.618'
$ws.Cells.Item(6, 5).Value = '  -0.25%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '74.97'
$ws.Cells.Item(7, 5).Value = '  +4.31%  '
$ws.Cells.Item(8, 5).Value = '  -0.05%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.603'
$ws.Cells.Item(9, 5).Value = '  +2.04%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '41.15'
$ws.Cells.Item(10, 5).Value = '  -0.75%  '
$ws.Cells.Item(11, 5).Value = '  -1.33%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '55.25'
$ws.Cells.Item(12, 5).Value = '  -1.35%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '6.89'
$ws.Cells.Item(13, 5).Value = '  +0.07%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.102'
$ws.Cells.Item(14, 5).Value = '  -0.90%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.552.11'
$ws.Cells.Item(15, 5).Value = '  +1.03%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.62'
$ws.Cells.Item(16, 5).Value = '  +3.18%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '2.223.37'
$ws.Cells.Item(17, 5).Value = '  +1.91%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.808'
$ws.Cells.Item(18, 5).Value = '  -2.20%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '42.613.16'
$ws.Cells.Item(19, 5).Value = '  +3.01%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0000104'
$ws.Cells.Item(20, 5).Value = '  -1.16%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '70.90'
$ws.Cells.Item(21, 5).Value = '  -0.77%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.93'
$ws.Cells.Item(22, 5).Value = '  -2.61%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '9.92'
$ws.Cells.Item(23, 5).Value = '  -4.51%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '229.60'
$ws.Cells.Item(24, 5).Value = '  +0.48%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.15'
$ws.Cells.Item(25, 5).Value = '  +7.04%  '
$ws.Cells.Item(26, 5).Value = '  -0.12%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.90'
$ws.Cells.Item(27, 5).Value = '  -2.55%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '3.33'
$ws.Cells.Item(28, 5).Value = '  -8.23%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.23'
$ws.Cells.Item(29, 5).Value = '  -0.94%  '
$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.19'
$ws.Cells.Item(30, 5).Value = '  -0.64%  '
$ws.Cells.Item(31, 2).Value = 'Monero'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '174.06'
$ws.Cells.Item(31, 5).Value = '  +4.51%  '
$ws.Cells.Item(32, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '36.85'
$ws.Cells.Item(32, 5).Value = '  +12.82%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '20.29'
$ws.Cells.Item(33, 5).Value = '  +0.12%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0791'
$ws.Cells.Item(34, 5).Value = '  -0.59%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.31'
$ws.Cells.Item(35, 5).Value = '  +0.93%  '
$ws.Cells.Item(36, 5).Value = '  -0.38%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.107'
$ws.Cells.Item(37, 5).Value = '  -0.25%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '4.39'
$ws.Cells.Item(38, 5).Value = '  +4.25%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0335'
$ws.Cells.Item(39, 5).Value = '  +13.22%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '12.99'
$ws.Cells.Item(40, 5).Value = '  +1.07%  '
$ws.Cells.Item(41, 5).Value = '  +1.48%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.53'
$ws.Cells.Item(42, 5).Value = '  -0.26%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '59.96'
$ws.Cells.Item(43, 5).Value = '  -2.92%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.197'
$ws.Cells.Item(44, 5).Value = '  +1.41%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '8.60'
$ws.Cells.Item(45, 5).Value = '  +0.87%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0991'
$ws.Cells.Item(46, 5).Value = '  +0.53%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '99.94'
$ws.Cells.Item(47, 5).Value = '  -0.95%  '
$ws.Cells.Item(48, 5).Value = '  -0.98%  '
$ws.Cells.Item(49, 2).Value = 'WOONetwork'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.436'
$ws.Cells.Item(49, 5).Value = '  +15.78%  '
$ws.Cells.Item(50, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.15'
$ws.Cells.Item(50, 5).Value = '  -0.73%  '
$ws.Cells.Item(51, 2).Value = 'NEARProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.28'
$ws.Cells.Item(51, 5).Value = '  -0.06%  '
